$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains a new "2020" data column (Q), to the right of the existing
# "2019" column (P). Start by cloning column P's cell formatting (borders,
# number format, font, alignment) into column Q for the header block and
# every data row so the new column visually matches the rest of the table.
$ws.Range("P3:P14").Copy() | Out-Null
$ws.Range("Q3:Q14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 4 header: the new year label.
$ws.Range("Q4").Value = 2020

# Rows 5-14: the 2020 figures for each region (row 3 / Q3 stays blank, it's
# just the thin separator row above the header).
$ws.Range("Q5").Value = 38.6
$ws.Range("Q6").Value = 42.4
$ws.Range("Q7").Value = 53.2
$ws.Range("Q8").Value = 90.6
$ws.Range("Q9").Value = 52.6
$ws.Range("Q10").Value = 24.5
$ws.Range("Q11").Value = 69.099999999999994
$ws.Range("Q12").Value = 32.200000000000003
$ws.Range("Q13").Value = 19.100000000000001
$ws.Range("Q14").Value = 25.2

# Update the recorded cursor position to match post-edit state.
$ws.Range("R27").Select() | Out-Null
